$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Work Area" row (row 10) -> "Tissue (General)"
$ws.Range("A10").Value = "Tissue (General)"
$ws.Range("C10").Value = 0.016
$ws.Range("F10").Value = 32

# Add new row 11: Skeletal Muscle
$ws.Range("A11").Value = "Skeletal Muscle"
$ws.Range("B11").Value = 45
$ws.Range("C11").Value = 0.01
$ws.Range("D11").Value = 0.036
$ws.Range("E11").Formula = "=(C11*(F11+B11))/(B11*D11)"
$ws.Range("F11").Value = 32
